$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista1")

# Fill in the previously empty row 4 (gap between the real data in rows 1-3
# and the zero-filled block that used to start at row 5) with zeros,
# matching the rest of the zero-filled block.
$ws.Range("A4:D4").Value = 0.0

# Remove the last row of the zero-filled block (row 12) so the total
# number of data rows stays the same as before.
$ws.Rows.Item(12).Delete()
